$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 163; existing rows 163:223 shift down to 164:224,
# matching the new sheet dimension A1:R224.
$ws.Rows("163:163").Insert()

# Populate the newly inserted row 163 with the new record.
$ws.Range("A163").Value = 5
$ws.Range("B163").Value = "Macroferia Regional de Talca"
$ws.Range("C163").Value = "Maule"
$ws.Range("D163").Value = "2021-09-29"
$ws.Range("E163").Value = 7
$ws.Range("F163").Value = 100112043
$ws.Range("G163").Value = "Pepino ensalada"
$ws.Range("H163").Value = "Sin especificar"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 300
$ws.Range("K163").Value = 16000
$ws.Range("L163").Value = 16000
$ws.Range("M163").Value = 16000
$ws.Range("N163").Value = "$/caja 60 unidades"
$ws.Range("O163").Value = "Región de Arica y Parinacota"
$ws.Range("P163").Value = 267
$ws.Range("Q163").Value = 60
$ws.Range("R163").Value = "Hortaliza"
